# first week updates to slides
#
# Slide 3 ("Examples from class"): add a new bullet
# "HTML Basics repo, Week CSS folder" before the existing
# "Margins and Padding" paragraph in the content placeholder.
#
# (Note: the source deck's slide master / layouts also carry an
# auto-updating "datetimeFigureOut" date field whose cached display
# text moved from 7/29/22 to 8/26/22 between commits. That value is
# refreshed by PowerPoint itself from the system clock whenever the
# field recalculates, rather than being literal content - it is left
# alone here so the <a:fld> element (id/type) stays intact instead of
# being collapsed into a plain run.)

$p = $ppt.ActivePresentation

$s3 = $p.Slides.Item(3)
for ($i = 1; $i -le $s3.Shapes.Count; $i++) {
    $shape = $s3.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        if ($shape.TextFrame.TextRange.Text -match "Margins and Padding") {
            $shape.TextFrame.TextRange.InsertBefore("HTML Basics repo, Week CSS folder`r") | Out-Null
        }
    }
}
